$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.353.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.75%  "

$ws.Range("D3").Value = "'3.442.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.56%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "'610.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.66%  "

$ws.Range("D6").Value = "'167.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.78%  "

$ws.Range("D7").Value = "'3.433.98"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.65%  "

$ws.Range("E8").Value = "  -2.15%  "

$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("D11").Value = "'7.08"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.04%  "

$ws.Range("E12").Value = "  -2.72%  "

$ws.Range("D13").Value = "'44.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.64%  "

$ws.Range("D14").Value = "'0.0000270"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.48%  "

$ws.Range("D15").Value = "'3.998.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.58%  "

$ws.Range("D16").Value = "'8.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.33%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'3.452.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.45%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'581.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.63%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "'69.369.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("D20").Value = "'0.121"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.11%  "

$ws.Range("D21").Value = "'17.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("D22").Value = "'0.846"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.71%  "

$ws.Range("D23").Value = "'8.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.49%  "

$ws.Range("D24").Value = "'95.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").Value = "'15.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.10%  "

$ws.Range("D26").Value = "'3.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.06%  "

$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("E28").Value = "  -5.03%  "

$ws.Range("D29").Value = "'32.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.04%  "

$ws.Range("D30").Value = "'8.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.28%  "

$ws.Range("D31").Value = "'7.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.20%  "

$ws.Range("E32").Value = "  -2.34%  "

$ws.Range("D33").Value = "'2.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.87%  "

$ws.Range("D34").Value = "'6.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.81%  "

$ws.Range("D35").Value = "'575.78"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.32%  "

$ws.Range("E36").Value = "  -1.49%  "

$ws.Range("D37").Value = "'0.0475"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").Value = "'0.0956"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.85%  "

$ws.Range("E39").Value = "  +0.17%  "

$ws.Range("D40").Value = "'55.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "

$ws.Range("D41").Value = "'0.141"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.82%  "

$ws.Range("D42").Value = "'3.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -11.51%  "

$ws.Range("D43").Value = "'3.242.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.25%  "

$ws.Range("D44").Value = "'0.0₃0684"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("D45").Value = "'31.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.24%  "

$ws.Range("D46").Value = "'0.295"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.68%  "

$ws.Range("E47").Value = "  -5.08%  "

$ws.Range("E48").Value = "  -5.50%  "

$ws.Range("E49").Value = "  -2.56%  "

$ws.Range("D50").Value = "'133.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("E51").Value = "  -0.03%  "
